$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and Report Covering the Week) ---
$a8 = $ws.Range("A8")
$a8.Value = "Volume 31   Number  27"
$a8.Characters(1,7).Font.Size = 10
$a8.Characters(1,7).Font.Name = "Andale WT"
$a8.Characters(8,2).Font.Size = 10
$a8.Characters(8,2).Font.Name = "Andale WT"
$a8.Characters(10,11).Font.Size = 10
$a8.Characters(10,11).Font.Name = "Andale WT"
$a8.Characters(21,2).Font.Size = 10
$a8.Characters(21,2).Font.Name = "Andale WT"

$c9 = $ws.Range("C9")
$c9.Value = "Report Covering the Week  7/1/2024  Through  7/7/2024"
$c9.Characters(1,26).Font.Size = 10
$c9.Characters(1,26).Font.Name = "Andale WT"
$c9.Characters(27,8).Font.Size = 10
$c9.Characters(27,8).Font.Name = "Andale WT"
$c9.Characters(35,11).Font.Size = 10
$c9.Characters(35,11).Font.Name = "Andale WT"
$c9.Characters(46,8).Font.Size = 10
$c9.Characters(46,8).Font.Name = "Andale WT"

# --- Cells changing style/type (use PasteSpecial from a donor cell holding the target style) ---
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C15").Value = 2

$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("K14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("C16").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("F15").Value = 2

$ws.Range("C14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("K14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("C16").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C27").Value = 2

$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("K14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("C16").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("F27").Value = 2

$ws.Range("C16").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("D28").Value = 1

$ws.Range("K22").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("E28").Value = -100

$ws.Range("C14").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163) | Out-Null  # values

$ws.Range("K14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null  # formats
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E31").PasteSpecial(-4163) | Out-Null  # values

$excel.CutCopyMode = $false

# --- Simple value-only updates (style unchanged) ---
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -14.285714285714
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -86.666666666666
$ws.Range("I16").Value = 38
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = -22.448979591836
$ws.Range("L16").Value = -9.523809523809
$ws.Range("M16").Value = 11.764705882352
$ws.Range("N16").Value = -87.936507936507
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 15.384615384615
$ws.Range("L17").Value = 3.448275862068
$ws.Range("M17").Value = 71.428571428571
$ws.Range("N17").Value = -36.842105263157
$ws.Range("C18").Value = 2
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 44
$ws.Range("K18").Value = -25.423728813559
$ws.Range("L18").Value = -58.878504672897
$ws.Range("M18").Value = -22.807017543859
$ws.Range("N18").Value = -92.334494773519
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 44
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -8.333333333333
$ws.Range("I19").Value = 296
$ws.Range("J19").Value = 378
$ws.Range("K19").Value = -21.693121693121
$ws.Range("L19").Value = -8.641975308641
$ws.Range("M19").Value = -18.232044198895
$ws.Range("N19").Value = -73.429084380610
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = -51.851851851851
$ws.Range("L20").Value = -55.172413793103
$ws.Range("N20").Value = -96.036585365853
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -45
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 81
$ws.Range("H21").Value = -24.691358024691
$ws.Range("I21").Value = 458
$ws.Range("J21").Value = 570
$ws.Range("K21").Value = -19.649122807017
$ws.Range("L21").Value = -19.081272084805
$ws.Range("M21").Value = -9.306930693069
$ws.Range("N21").Value = -81.175503493629
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 33.333333333333
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = -55.263157894736
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 111
$ws.Range("H24").Value = -29.729729729729
$ws.Range("I24").Value = 497
$ws.Range("J24").Value = 526
$ws.Range("K24").Value = -5.513307984790
$ws.Range("L24").Value = -23.655913978494
$ws.Range("M24").Value = 51.063829787234
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 38
$ws.Range("E25").Value = -71.052631578947
$ws.Range("G25").Value = 93
$ws.Range("H25").Value = -32.258064516129
$ws.Range("I25").Value = 404
$ws.Range("J25").Value = 433
$ws.Range("K25").Value = -6.697459584295
$ws.Range("L25").Value = -25.598526703499
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 15
$ws.Range("H26").Value = -13.333333333333
$ws.Range("I26").Value = 104
$ws.Range("J26").Value = 138
$ws.Range("K26").Value = -24.637681159420
$ws.Range("L26").Value = -13.333333333333
$ws.Range("M26").Value = -18.110236220472
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -12.5
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 29.166666666666
